$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "formulation" column in K
$ws.Range("K1").Value = "formulation"
$ws.Range("K2").Value = "ADS-5102 capsule"
$ws.Range("K3").Value = "Amantadine hydrochloride"
$ws.Range("K4").Value = "ADS-5102 capsule"
$ws.Range("K5").Value = "ADS-5102 capsule"
$ws.Range("K6").Value = "Amantadine hydrochloride"
$ws.Range("K7").Value = "Other/not specified"
$ws.Range("K8").Value = "Other/not specified"
$ws.Range("K9").Value = "Amantadine sulfate"
$ws.Range("K10").Value = "ADS-5102 capsule"

# Update the selection to match the saved view state in the diff
$ws.Range("J13").Select()
